$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill column E (round 1 results) with the winner picked in column D for most rows.
# Rows 2,4,6,8,10,12,16 -> same winner as column D (copy across)
$ws.Range("E2").Value2 = $ws.Range("D2").Value2
$ws.Range("E4").Value2 = $ws.Range("D4").Value2
$ws.Range("E6").Value2 = $ws.Range("D6").Value2
$ws.Range("E8").Value2 = $ws.Range("D8").Value2
$ws.Range("E10").Value2 = $ws.Range("D10").Value2
$ws.Range("E12").Value2 = $ws.Range("D12").Value2
$ws.Range("E16").Value2 = $ws.Range("D16").Value2

# Row 14 is different: fixed bracket so the winner shown in E14 is the opponent from D15
$ws.Range("E14").Value2 = $ws.Range("D15").Value2

# Update the selection / view state to match the saved workbook
$ws.Range("E4").Select()

$wb.Save()
